$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 322.2857
$ws.Range("I2").Value = 302
$ws.Range("K2").Value = 302
$ws.Range("M2").Value = -189

$ws.Range("H18").Value = 224.21739
$ws.Range("I18").Value = 211.68182
$ws.Range("K18").Value = 211.68182
$ws.Range("M18").Value = 72.31818000000001

$ws.Range("H40").Value = 918.6667
$ws.Range("I40").Value = 882.8570999999999
$ws.Range("J40").Value = 950
$ws.Range("K40").Value = 882.8570999999999
$ws.Range("L40").Value = 950
$ws.Range("M40").Value = -707.8570999999999
$ws.Range("N40").Value = -1300

$ws.Range("H76").Value = 66669616
$ws.Range("I76").Value = 66669616
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 66669616
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -66669301
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 66669616
$ws.Range("I79").Value = 66669616
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 66669616
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -66668524
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 4397.25
$ws.Range("I5").Value = 5699.6665
$ws.Range("J5").Value = 490
$ws.Range("K5").Value = 5699.6665
$ws.Range("L5").Value = 490
$ws.Range("M5").Value = -5587.6665
$ws.Range("N5").Value = -714

$ws.Range("H63").Value = 2017
$ws.Range("I63").Value = 2017
$ws.Range("K63").Value = 2017
$ws.Range("M63").Value = -1331

$ws.Range("H66").Value = 2017
$ws.Range("I66").Value = 2017
$ws.Range("K66").Value = 10085
$ws.Range("M66").Value = -6653

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 4397.25
$ws.Range("I4").Value = 5699.6665
$ws.Range("J4").Value = 490
$ws.Range("K4").Value = 5699.6665
$ws.Range("L4").Value = 490
$ws.Range("M4").Value = -5584.6665
$ws.Range("N4").Value = -720

$ws.Range("H22").Value = 494.29413
$ws.Range("I22").Value = 450.5
$ws.Range("J22").Value = 500.13333
$ws.Range("K22").Value = 450.5
$ws.Range("L22").Value = 500.13333
$ws.Range("M22").Value = -277.5
$ws.Range("N22").Value = -846.13333

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H82").Value = 5680
$ws.Range("I82").Value = 5680
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 5680
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -5297
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 5680
$ws.Range("I85").Value = 5680
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 5680
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -4354
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 734.05
$ws.Range("I22").Value = 814.9375
$ws.Range("K22").Value = 814.9375
$ws.Range("M22").Value = -464.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 901
$ws.Range("I92").Value = 802
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 2406
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -1158
$ws.Range("N92").Value = -5496

$ws.Range("H97").Value = 8252.5
$ws.Range("J97").Value = 8252.5
$ws.Range("L97").Value = 24757.5
$ws.Range("N97").Value = -25749.5

$ws.Range("H115").Value = 4000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 4000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 12000
$ws.Range("N115").Value = -14350
$ws.Range("M115").ClearContents()

$ws.Range("H131").Value = 40599210
$ws.Range("I131").Value = 846.5
$ws.Range("J131").Value = 47980730
$ws.Range("K131").Value = 2539.5
$ws.Range("L131").Value = 143942190
$ws.Range("M131").Value = 2500.5
$ws.Range("N131").Value = -143952270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H43").Value = 9615.223
$ws.Range("I43").Value = 700
$ws.Range("J43").Value = 12162.429
$ws.Range("K43").Value = 700
$ws.Range("L43").Value = 12162.429
$ws.Range("M43").Value = -549
$ws.Range("N43").Value = -12464.429

$ws.Range("H55").Value = 62516.5
$ws.Range("J55").Value = 62516.5
$ws.Range("L55").Value = 62516.5
$ws.Range("N55").Value = -63170.5

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 10000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 10000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -9617
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 10000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 10000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -8674
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 585.5263
$ws.Range("I22").Value = 502.5
$ws.Range("J22").Value = 677.7778
$ws.Range("K22").Value = 502.5
$ws.Range("L22").Value = 677.7778
$ws.Range("M22").Value = -207.5
$ws.Range("N22").Value = -1267.7778

$ws.Range("H27").Value = 585.5263
$ws.Range("I27").Value = 502.5
$ws.Range("J27").Value = 677.7778
$ws.Range("K27").Value = 502.5
$ws.Range("L27").Value = 677.7778
$ws.Range("M27").Value = -395.5
$ws.Range("N27").Value = -891.7778

$ws.Range("H95").Value = 7724.75
$ws.Range("J95").Value = 7724.75
$ws.Range("L95").Value = 7724.75
$ws.Range("N95").Value = -13216.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 26091.666
$ws.Range("J97").Value = 26091.666
$ws.Range("L97").Value = 26091.666
$ws.Range("N97").Value = -28073.666

